# Build site at 2022-09-26 16:07:08 UTC
#
# The LOM3095 course-info sheet had several label/value rows reshuffled:
#  - "Objetivos:" (row 10) now shows the "5840521 - Rosa Ana Conte" value
#  - a new "Programa resumido:" / "Semestral" pair lands on row 13
#  - the long free-text paragraphs that used to sit on rows 14 and 16 are
#    dropped
#  - "Programa:" (row 15) now shows the "01/01/2020" value
#  - "Método:" (row 18) now shows the "5840521 - Rosa Ana Conte" value
#  - the row-13..21 labels shift down by one (Short syllabus/Programa/
#    Syllabus/Avaliação/Método/Critério/Norma de recuperação/Bibliografia)
#  - the old last row (22, "Bibliografia:" + bibliography text) disappears,
#    shrinking the sheet from A1:C22 to A1:C21
#
# Note: new B/C cells are populated via Range.Copy() from a same-column
# cell that already carries the right text/style, rather than via
# Range.Value =, to avoid two COM quirks: (1) a literal "01/01/2020" typed
# into .Value gets reinterpreted as a date serial instead of staying the
# original text, and (2) this sheet's <cols> has an overlapping style
# definition for column B, so a brand-new B-column cell created via
# .Value = picks up the row's style instead of the column's.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the old trailing row entirely (label + bibliography text) ---
$ws.Rows(22).Delete()

# --- Row 10: value changes to the professor string (label stays) ---
$ws.Range("B10").Value = "5840521 - Rosa Ana Conte"
$ws.Range("C10").Value = "5840521 - Rosa Ana Conte"

# --- Row 18 will need the same "5840521 - Rosa Ana Conte" text; copy it
#     now from B13/C13 (which still hold it) before row 13 is overwritten
#     below, so the new B18/C18 cells inherit the correct column style. ---
$ws.Range("B13").Copy($ws.Range("B18"))
$ws.Range("C13").Copy($ws.Range("C18"))

# --- Row 13 gains a label + new "Semestral" value ---
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# --- Row 14 keeps only its (new) label; old paragraph text is dropped ---
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14:C14").Clear()

# --- Row 15 gets its (new) label and the "01/01/2020" value; copy the
#     existing text cell from B8/C8 so it stays a literal text value
#     instead of being reinterpreted as a date, and picks up the right
#     column style as a bonus. ---
$ws.Range("A15").Value = "Programa:"
$ws.Range("B8").Copy($ws.Range("B15"))
$ws.Range("C8").Copy($ws.Range("C15"))

# --- Row 16 keeps only its (new) label; old paragraph text is dropped ---
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16:C16").Clear()

# --- Remaining label shifts (rows 17-21); B/C values already correct ---
$ws.Range("A17").Value = "Avaliação:"
$ws.Range("A18").Value = "Método:"
$ws.Range("A19").Value = "Critério:"
$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("A21").Value = "Bibliografia:"

# --- Row heights: a few rows gain/change a custom height, one loses it ---
$ws.Rows(13).RowHeight = 60
$ws.Rows(15).RowHeight = 120
$ws.Rows(17).AutoFit()
$ws.Rows(18).RowHeight = 60
$ws.Rows(21).RowHeight = 120
